$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "Data2"
$ws.Range("G2").Value = 5
$ws.Range("G3").Value = 6
$ws.Range("G4").Value = 7
$ws.Range("G5").Value = 89

$ws.Range("G6").Select()
